$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, pushing the existing rows 10-18 down to 11-19.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly price record.
$ws.Cells.Item(10, 1).Value = 11
$ws.Cells.Item(10, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(10, 3).Value = "Bíobío"
$ws.Cells.Item(10, 4).Value = 45100
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100107
$ws.Cells.Item(10, 8).Value = "Otros"
$ws.Cells.Item(10, 9).Value = 100107001
$ws.Cells.Item(10, 10).Value = "Caqui"
$ws.Cells.Item(10, 11).Value = "Mankaki"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 60
$ws.Cells.Item(10, 14).Value = 18000
$ws.Cells.Item(10, 15).Value = 18000
$ws.Cells.Item(10, 16).Value = 18000
$ws.Cells.Item(10, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(10, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 19).Value = 1000
$ws.Cells.Item(10, 20).Value = 18
